# n-dimensional euclidean example
# Adds a new "Sheet2" worksheet (K-Nearest Neighbor, euclidean distance over
# weight/height) after the existing sheets, populates it with the sample
# data + SQRT/SMALL formulas, formats the header row + distance columns,
# adds a threaded comment with a reference link, and updates the
# selection/active-sheet state to match.

$wb = $excel.ActiveWorkbook

# --- Move the current selection on "K-Nearest Neighbor" before adding the
#     new sheet, so it stops being the tab-selected sheet and its cursor
#     moves from A12 to F5.
$knn = $wb.Worksheets.Item("K-Nearest Neighbor")
$knn.Activate()
$knn.Range("F5").Select()

# --- Add the new worksheet as the last tab, named "Sheet2".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sheet2"

# --- Sample data: weight (A), height (B), class (C), distance formula (D)
#     reference point weight/height (F2/G2), and nearest-3 distances (F4:F6).
# String literals are written in this specific order so the shared-string
# table is built up the same way the source workbook has it.
$ws.Range("C1").Value = "class"

$ws.Range("A2").Value = 51
$ws.Range("B2").Value = 167
$ws.Range("C2").Value = "underweight"

$ws.Range("A3").Value = 62
$ws.Range("B3").Value = 182
$ws.Range("C3").Value = "normal"

$ws.Range("A4").Value = 69
$ws.Range("B4").Value = 176
$ws.Range("C4").Value = "normal"

$ws.Range("A5").Value = 64
$ws.Range("B5").Value = 173
$ws.Range("C5").Value = "normal"

$ws.Range("A6").Value = 65
$ws.Range("B6").Value = 172
$ws.Range("C6").Value = "normal"

$ws.Range("A7").Value = 56
$ws.Range("B7").Value = 174
$ws.Range("C7").Value = "underweight"

$ws.Range("A8").Value = 58
$ws.Range("B8").Value = 169
$ws.Range("C8").Value = "normal"

$ws.Range("A9").Value = 57
$ws.Range("B9").Value = 173
$ws.Range("C9").Value = "normal"

$ws.Range("A10").Value = 55
$ws.Range("B10").Value = 170
$ws.Range("C10").Value = "normal"

$ws.Range("D1").Value = "distance"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "height"

$ws.Range("F2").Value = 57
$ws.Range("G2").Value = 170

$ws.Range("F3").Value = "k=3"

# --- Euclidean distance of each row's (weight,height) from the reference
#     point held in $F$2/$G$2. Single-quoted so PowerShell does not try to
#     interpolate the $G$2 / $F$2 / $D$2 style absolute references.
$ws.Range("D2").Formula = '=SQRT(($G$2-B2) * ($G$2-B2) + ($F$2-A2) * ($F$2-A2))'
$ws.Range("D3").Formula = '=SQRT(($G$2-B3) * ($G$2-B3) + ($F$2-A3) * ($F$2-A3))'
$ws.Range("D4").Formula = '=SQRT(($G$2-B4) * ($G$2-B4) + ($F$2-A4) * ($F$2-A4))'
$ws.Range("D5").Formula = '=SQRT(($G$2-B5) * ($G$2-B5) + ($F$2-A5) * ($F$2-A5))'
$ws.Range("D6").Formula = '=SQRT(($G$2-B6) * ($G$2-B6) + ($F$2-A6) * ($F$2-A6))'
$ws.Range("D7").Formula = '=SQRT(($G$2-B7) * ($G$2-B7) + ($F$2-A7) * ($F$2-A7))'
$ws.Range("D8").Formula = '=SQRT(($G$2-B8) * ($G$2-B8) + ($F$2-A8) * ($F$2-A8))'
$ws.Range("D9").Formula = '=SQRT(($G$2-B9) * ($G$2-B9) + ($F$2-A9) * ($F$2-A9))'
$ws.Range("D10").Formula = '=SQRT(($G$2-B10) * ($G$2-B10) + ($F$2-A10) * ($F$2-A10))'

# --- 3 smallest distances (k=3 nearest neighbours).
$ws.Range("F4").Formula = '=SMALL($D$2:$D$10,ROWS(C$2:C2))'
$ws.Range("F5").Formula = '=SMALL($D$2:$D$10,ROWS(C$2:C3))'
$ws.Range("F6").Formula = '=SMALL($D$2:$D$10,ROWS(C$2:C4))'

# --- Number formatting: one decimal place on the distance/SMALL columns.
# Applied before the bold header font so the numFmt style is created first
# (matches the style index ordering of the source file).
$ws.Range("D2:D10").NumberFormat = "0.0"
$ws.Range("F4:F6").NumberFormat = "0.0"

# --- Bold header row + the "k=3" label.
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F3").Font.Bold = $true

# --- Column widths roughly matching the source (best-fit to content).
# The host's ColumnWidth setter stores (input + 0.8333) "characters", so the
# inputs below are back-solved to land on the source file's stored widths.
$ws.Columns.Item(1).ColumnWidth = 5.830729166666667
$ws.Columns.Item(2).ColumnWidth = 5.498697916666667
$ws.Columns.Item(3).ColumnWidth = 10.276041666666666
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 3.0533854166666665

# --- Threaded comment on the "class" header, linking to the reference video.
$excel.UserName = "Randy Hollines"
$ws.Range("C1").AddCommentThreaded("https://www.youtube.com/watch?v=4HKqjENq9OU")

# --- Landscape/portrait page setup so a <pageSetup> element is emitted.
$ws.PageSetup.Orientation = 1

# --- Make the new sheet the active tab, cursor on F7 (matches the diff).
$ws.Activate()
$ws.Range("F7").Select()
